$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new first row (pushes everything else down by one) and
# give it the new "Hello World" command text that was missing from the
# output-parsing test fixture.
$ws.Rows("1").Insert()
$ws.Range("A1").Value = "$> echo 'Hello World' | grep Hello | .wc -w >> output.txt"

# The remaining rows (now rows 4 through 11, originally rows 3 through 10)
# exercised multi-line / ANSI-color-escape cases that are no longer needed,
# so drop them - only the first three rows survive.
$ws.Rows("4:11").Delete()

# Row 3 (previously styled with wrap-text + a taller row height for the
# multi-line samples) now holds a short single-line string, so restore the
# default style/height.
$ws.Range("A3").Style = "Normal"
$ws.Rows("3").AutoFit()

# Match the saved selection/active cell state.
$ws.Range("A3").Select()
